# Update cryptocurrency market data table (A2:G51) with refreshed values
# as of 2024-03-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rank 1: BTC (Bitcoin)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = 69792
$ws.Cells.Item(2, 5).Value = 1372768066343
$ws.Cells.Item(2, 6).Value = 29480287500
$ws.Cells.Item(2, 7).Value = 2.11147

# Rank 2: ETH (Ethereum)
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = 3947.06
$ws.Cells.Item(3, 5).Value = 474150365225
$ws.Cells.Item(3, 6).Value = 14032475055
$ws.Cells.Item(3, 7).Value = 0.87414

# Rank 3: USDT (Tether)
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = "Tether"
$ws.Cells.Item(4, 4).Value = 1.002
$ws.Cells.Item(4, 5).Value = 102068929003
$ws.Cells.Item(4, 6).Value = 57138276325
$ws.Cells.Item(4, 7).Value = 0.07492

# Rank 4: BNB (BNB)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "BNB"
$ws.Cells.Item(5, 4).Value = 535.26
$ws.Cells.Item(5, 5).Value = 82447329943
$ws.Cells.Item(5, 6).Value = 3733362208
$ws.Cells.Item(5, 7).Value = 9.94711

# Rank 5: SOL (Solana)
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "SOL"
$ws.Cells.Item(6, 3).Value = "Solana"
$ws.Cells.Item(6, 4).Value = 145.95
$ws.Cells.Item(6, 5).Value = 64706400180
$ws.Cells.Item(6, 6).Value = 3098432302
$ws.Cells.Item(6, 7).Value = -0.21787

# Rank 6: STETH (Lido Staked Ether)
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "STETH"
$ws.Cells.Item(7, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(7, 4).Value = 3941.66
$ws.Cells.Item(7, 5).Value = 38901061402
$ws.Cells.Item(7, 6).Value = 20746321
$ws.Cells.Item(7, 7).Value = 0.78876

# Rank 7: XRP (XRP)
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "XRP"
$ws.Cells.Item(8, 3).Value = "XRP"
$ws.Cells.Item(8, 4).Value = 0.621121
$ws.Cells.Item(8, 5).Value = 34008213135
$ws.Cells.Item(8, 6).Value = 1361061218
$ws.Cells.Item(8, 7).Value = -0.20643

# Rank 8: USDC (USDC)
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "USDC"
$ws.Cells.Item(9, 3).Value = "USDC"
$ws.Cells.Item(9, 4).Value = 0.999821
$ws.Cells.Item(9, 5).Value = 30163413511
$ws.Cells.Item(9, 6).Value = 6228037557
$ws.Cells.Item(9, 7).Value = 0.07329

# Rank 9: ADA (Cardano)
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "ADA"
$ws.Cells.Item(10, 3).Value = "Cardano"
$ws.Cells.Item(10, 4).Value = 0.730597
$ws.Cells.Item(10, 5).Value = 25703642667
$ws.Cells.Item(10, 6).Value = 564037720
$ws.Cells.Item(10, 7).Value = -0.02436

# Rank 10: DOGE (Dogecoin)
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "DOGE"
$ws.Cells.Item(11, 3).Value = "Dogecoin"
$ws.Cells.Item(11, 4).Value = 0.174068
$ws.Cells.Item(11, 5).Value = 24952662017
$ws.Cells.Item(11, 6).Value = 3529979595
$ws.Cells.Item(11, 7).Value = 4.90671

# Rank 11: SHIB (Shiba Inu)
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "SHIB"
$ws.Cells.Item(12, 3).Value = "Shiba Inu"
$ws.Cells.Item(12, 4).Value = 0.00003389
$ws.Cells.Item(12, 5).Value = 19980795015
$ws.Cells.Item(12, 6).Value = 2159557539
$ws.Cells.Item(12, 7).Value = -0.98049

# Rank 12: AVAX (Avalanche)
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "AVAX"
$ws.Cells.Item(13, 3).Value = "Avalanche"
$ws.Cells.Item(13, 4).Value = 42.84
$ws.Cells.Item(13, 5).Value = 16170083382
$ws.Cells.Item(13, 6).Value = 551695626
$ws.Cells.Item(13, 7).Value = -1.01987

# Rank 13: DOT (Polkadot)
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "DOT"
$ws.Cells.Item(14, 3).Value = "Polkadot"
$ws.Cells.Item(14, 4).Value = 10.44
$ws.Cells.Item(14, 5).Value = 13999088483
$ws.Cells.Item(14, 6).Value = 392533419
$ws.Cells.Item(14, 7).Value = -3.95771

# Rank 14: TRX (TRON)
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "TRX"
$ws.Cells.Item(15, 3).Value = "TRON"
$ws.Cells.Item(15, 4).Value = 0.13586
$ws.Cells.Item(15, 5).Value = 11947735137
$ws.Cells.Item(15, 6).Value = 349480182
$ws.Cells.Item(15, 7).Value = -0.1146

# Rank 15: LINK (Chainlink)
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "LINK"
$ws.Cells.Item(16, 3).Value = "Chainlink"
$ws.Cells.Item(16, 4).Value = 19.89
$ws.Cells.Item(16, 5).Value = 11689227698
$ws.Cells.Item(16, 6).Value = 464409434
$ws.Cells.Item(16, 7).Value = -0.67979

# Rank 16: MATIC (Polygon)
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "MATIC"
$ws.Cells.Item(17, 3).Value = "Polygon"
$ws.Cells.Item(17, 4).Value = 1.22
$ws.Cells.Item(17, 5).Value = 11349289853
$ws.Cells.Item(17, 6).Value = 896885173
$ws.Cells.Item(17, 7).Value = 7.41905

# Rank 17: WBTC (Wrapped Bitcoin)
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "WBTC"
$ws.Cells.Item(18, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(18, 4).Value = 69662
$ws.Cells.Item(18, 5).Value = 10881681787
$ws.Cells.Item(18, 6).Value = 201038582
$ws.Cells.Item(18, 7).Value = 1.78569

# Rank 18: UNI (Uniswap)
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "UNI"
$ws.Cells.Item(19, 3).Value = "Uniswap"
$ws.Cells.Item(19, 4).Value = 14.14
$ws.Cells.Item(19, 5).Value = 10635468364
$ws.Cells.Item(19, 6).Value = 302493958
$ws.Cells.Item(19, 7).Value = -0.73521

# Rank 19: TON (Toncoin)
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "TON"
$ws.Cells.Item(20, 3).Value = "Toncoin"
$ws.Cells.Item(20, 4).Value = 2.84
$ws.Cells.Item(20, 5).Value = 9868879057
$ws.Cells.Item(20, 6).Value = 64675832
$ws.Cells.Item(20, 7).Value = -1.26682

# Rank 20: BCH (Bitcoin Cash)
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "BCH"
$ws.Cells.Item(21, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(21, 4).Value = 434.32
$ws.Cells.Item(21, 5).Value = 8546928858
$ws.Cells.Item(21, 6).Value = 394775514
$ws.Cells.Item(21, 7).Value = 0.44527

# Rank 21: ICP (Internet Computer)
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "ICP"
$ws.Cells.Item(22, 3).Value = "Internet Computer"
$ws.Cells.Item(22, 4).Value = 14.6
$ws.Cells.Item(22, 5).Value = 6717111393
$ws.Cells.Item(22, 6).Value = 181197216
$ws.Cells.Item(22, 7).Value = -2.31756

# Rank 22: LTC (Litecoin)
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "LTC"
$ws.Cells.Item(23, 3).Value = "Litecoin"
$ws.Cells.Item(23, 4).Value = 88.75
$ws.Cells.Item(23, 5).Value = 6600340584
$ws.Cells.Item(23, 6).Value = 535070214
$ws.Cells.Item(23, 7).Value = 0.81255

# Rank 23: NEAR (NEAR Protocol)
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "NEAR"
$ws.Cells.Item(24, 3).Value = "NEAR Protocol"
$ws.Cells.Item(24, 4).Value = 6.09
$ws.Cells.Item(24, 5).Value = 6353400075
$ws.Cells.Item(24, 6).Value = 510322927
$ws.Cells.Item(24, 7).Value = -3.21122

# Rank 24: FIL (Filecoin)
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "FIL"
$ws.Cells.Item(25, 3).Value = "Filecoin"
$ws.Cells.Item(25, 4).Value = 10.88
$ws.Cells.Item(25, 5).Value = 5690830573
$ws.Cells.Item(25, 6).Value = 487998958
$ws.Cells.Item(25, 7).Value = -3.57936

# Rank 25: ETC (Ethereum Classic)
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "ETC"
$ws.Cells.Item(26, 3).Value = "Ethereum Classic"
$ws.Cells.Item(26, 4).Value = 36.79
$ws.Cells.Item(26, 5).Value = 5370960040
$ws.Cells.Item(26, 6).Value = 261645856
$ws.Cells.Item(26, 7).Value = -3.17887

# Rank 26: LEO (LEO Token)
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "LEO"
$ws.Cells.Item(27, 3).Value = "LEO Token"
$ws.Cells.Item(27, 4).Value = 5.7
$ws.Cells.Item(27, 5).Value = 5291579047
$ws.Cells.Item(27, 6).Value = 2616817
$ws.Cells.Item(27, 7).Value = -0.10607

# Rank 27: ATOM (Cosmos Hub)
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "ATOM"
$ws.Cells.Item(28, 3).Value = "Cosmos Hub"
$ws.Cells.Item(28, 4).Value = 13.34
$ws.Cells.Item(28, 5).Value = 5202848287
$ws.Cells.Item(28, 6).Value = 286292175
$ws.Cells.Item(28, 7).Value = -2.99241

# Rank 28: APT (Aptos)
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "APT"
$ws.Cells.Item(29, 3).Value = "Aptos"
$ws.Cells.Item(29, 4).Value = 13.07
$ws.Cells.Item(29, 5).Value = 4825470332
$ws.Cells.Item(29, 6).Value = 167653450
$ws.Cells.Item(29, 7).Value = -3.01964

# Rank 29: IMX (Immutable)
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "IMX"
$ws.Cells.Item(30, 3).Value = "Immutable"
$ws.Cells.Item(30, 4).Value = 3.39
$ws.Cells.Item(30, 5).Value = 4709282684
$ws.Cells.Item(30, 6).Value = 127754714
$ws.Cells.Item(30, 7).Value = -3.29582

# Rank 30: OP (Optimism)
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "OP"
$ws.Cells.Item(31, 3).Value = "Optimism"
$ws.Cells.Item(31, 4).Value = 4.65
$ws.Cells.Item(31, 5).Value = 4684388479
$ws.Cells.Item(31, 6).Value = 335141421
$ws.Cells.Item(31, 7).Value = 2.80675

# Rank 31: STX (Stacks)
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "STX"
$ws.Cells.Item(32, 3).Value = "Stacks"
$ws.Cells.Item(32, 4).Value = 3.24
$ws.Cells.Item(32, 5).Value = 4653182074
$ws.Cells.Item(32, 6).Value = 191954180
$ws.Cells.Item(32, 7).Value = 14.64591

# Rank 33: TAO (Bittensor)
$ws.Cells.Item(33, 1).Value = 33
$ws.Cells.Item(33, 2).Value = "TAO"
$ws.Cells.Item(33, 3).Value = "Bittensor"
$ws.Cells.Item(33, 4).Value = 705.29
$ws.Cells.Item(33, 5).Value = 4503492201
$ws.Cells.Item(33, 6).Value = 21293526
$ws.Cells.Item(33, 7).Value = -1.49537

# Rank 32: RNDR (Render)
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = "RNDR"
$ws.Cells.Item(34, 3).Value = "Render"
$ws.Cells.Item(34, 4).Value = 11.8
$ws.Cells.Item(34, 5).Value = 4497026052
$ws.Cells.Item(34, 6).Value = 651198464
$ws.Cells.Item(34, 7).Value = 5.29091

# Rank 34: DAI (Dai)
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "DAI"
$ws.Cells.Item(35, 3).Value = "Dai"
$ws.Cells.Item(35, 4).Value = 1.002
$ws.Cells.Item(35, 5).Value = 4459558988
$ws.Cells.Item(35, 6).Value = 241188651
$ws.Cells.Item(35, 7).Value = 0.11011

# Rank 35: CRO (Cronos)
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "CRO"
$ws.Cells.Item(36, 3).Value = "Cronos"
$ws.Cells.Item(36, 4).Value = 0.167438
$ws.Cells.Item(36, 5).Value = 4445744063
$ws.Cells.Item(36, 6).Value = 57181585
$ws.Cells.Item(36, 7).Value = 3.5072

# Rank 36: GRT (The Graph)
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "GRT"
$ws.Cells.Item(37, 3).Value = "The Graph"
$ws.Cells.Item(37, 4).Value = 0.471982
$ws.Cells.Item(37, 5).Value = 4410355072
$ws.Cells.Item(37, 6).Value = 1157517684
$ws.Cells.Item(37, 7).Value = 20.43049

# Rank 37: HBAR (Hedera)
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "HBAR"
$ws.Cells.Item(38, 3).Value = "Hedera"
$ws.Cells.Item(38, 4).Value = 0.128049
$ws.Cells.Item(38, 5).Value = 4310912076
$ws.Cells.Item(38, 6).Value = 73505372
$ws.Cells.Item(38, 7).Value = -1.23158

# Rank 38: OKB (OKB)
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "OKB"
$ws.Cells.Item(39, 3).Value = "OKB"
$ws.Cells.Item(39, 4).Value = 71.33
$ws.Cells.Item(39, 5).Value = 4299541978
$ws.Cells.Item(39, 6).Value = 53305462
$ws.Cells.Item(39, 7).Value = 17.98066

# Rank 39: XLM (Stellar)
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "XLM"
$ws.Cells.Item(40, 3).Value = "Stellar"
$ws.Cells.Item(40, 4).Value = 0.14272
$ws.Cells.Item(40, 5).Value = 4087486935
$ws.Cells.Item(40, 6).Value = 124480744
$ws.Cells.Item(40, 7).Value = 1.0794

# Rank 40: PEPE (Pepe)
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "PEPE"
$ws.Cells.Item(41, 3).Value = "Pepe"
$ws.Cells.Item(41, 4).Value = 0.00000854
$ws.Cells.Item(41, 5).Value = 3605453006
$ws.Cells.Item(41, 6).Value = 1262647504
$ws.Cells.Item(41, 7).Value = -2.15797

# Rank 41: INJ (Injective)
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "INJ"
$ws.Cells.Item(42, 3).Value = "Injective"
$ws.Cells.Item(42, 4).Value = 40.59
$ws.Cells.Item(42, 5).Value = 3593120848
$ws.Cells.Item(42, 6).Value = 156925788
$ws.Cells.Item(42, 7).Value = -1.94281

# Rank 42: VET (VeChain)
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "VET"
$ws.Cells.Item(43, 3).Value = "VeChain"
$ws.Cells.Item(43, 4).Value = 0.04839979
$ws.Cells.Item(43, 5).Value = 3520741472
$ws.Cells.Item(43, 6).Value = 136526334
$ws.Cells.Item(43, 7).Value = -0.35568

# Rank 43: MNT (Mantle)
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "MNT"
$ws.Cells.Item(44, 3).Value = "Mantle"
$ws.Cells.Item(44, 4).Value = 1.062
$ws.Cells.Item(44, 5).Value = 3433028163
$ws.Cells.Item(44, 6).Value = 218531795
$ws.Cells.Item(44, 7).Value = 2.19989

# Rank 44: KAS (Kaspa)
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "KAS"
$ws.Cells.Item(45, 3).Value = "Kaspa"
$ws.Cells.Item(45, 4).Value = 0.149773
$ws.Cells.Item(45, 5).Value = 3428401596
$ws.Cells.Item(45, 6).Value = 83972213
$ws.Cells.Item(45, 7).Value = 2.575

# Rank 45: FDUSD (First Digital USD)
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "FDUSD"
$ws.Cells.Item(46, 3).Value = "First Digital USD"
$ws.Cells.Item(46, 4).Value = 1.003
$ws.Cells.Item(46, 5).Value = 3303587057
$ws.Cells.Item(46, 6).Value = 9401755959
$ws.Cells.Item(46, 7).Value = -0.06571

# Rank 46: ARB (Arbitrum)
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "ARB"
$ws.Cells.Item(47, 3).Value = "Arbitrum"
$ws.Cells.Item(47, 4).Value = 2.1
$ws.Cells.Item(47, 5).Value = 3042844804
$ws.Cells.Item(47, 6).Value = 520026168
$ws.Cells.Item(47, 7).Value = -1.68752

# Rank 47: THETA (Theta Network)
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "THETA"
$ws.Cells.Item(48, 3).Value = "Theta Network"
$ws.Cells.Item(48, 4).Value = 3.01
$ws.Cells.Item(48, 5).Value = 3012862586
$ws.Cells.Item(48, 6).Value = 72720603
$ws.Cells.Item(48, 7).Value = -3.78113

# Rank 48: LDO (Lido DAO)
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "LDO"
$ws.Cells.Item(49, 3).Value = "Lido DAO"
$ws.Cells.Item(49, 4).Value = 3.34
$ws.Cells.Item(49, 5).Value = 2977986124
$ws.Cells.Item(49, 6).Value = 82529926
$ws.Cells.Item(49, 7).Value = -1.67793

# Rank 49: FET (Fetch.ai)
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "FET"
$ws.Cells.Item(50, 3).Value = "Fetch.ai"
$ws.Cells.Item(50, 4).Value = 2.81
$ws.Cells.Item(50, 5).Value = 2937498258
$ws.Cells.Item(50, 6).Value = 627542252
$ws.Cells.Item(50, 7).Value = -2.79714

# Rank 50: TIA (Celestia)
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "TIA"
$ws.Cells.Item(51, 3).Value = "Celestia"
$ws.Cells.Item(51, 4).Value = 16.46
$ws.Cells.Item(51, 5).Value = 2789998305
$ws.Cells.Item(51, 6).Value = 135259188
$ws.Cells.Item(51, 7).Value = -4.48562

Write-Output "Updated 50 rows (A2:G51)"
